$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. It belongs right
# before the existing row 37 (chronologically it is the most recent
# entry), so insert a fresh row there and push the rest of the table
# down by one -- this matches rows 37-40 (old) becoming rows 38-41 (new).
$ws.Rows.Item(37).Insert()

$newRow = 37

$ws.Cells.Item($newRow, 1).Value = 3
$ws.Cells.Item($newRow, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value = 44463
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat
$ws.Cells.Item($newRow, 5).Value = 5
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100107
$ws.Cells.Item($newRow, 8).Value = "Otros"
$ws.Cells.Item($newRow, 9).Value = 100107002
$ws.Cells.Item($newRow, 10).Value = "Chirimoya"
$ws.Cells.Item($newRow, 11).Value = "Cultivar IV Región"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 45
$ws.Cells.Item($newRow, 14).Value = 26000
$ws.Cells.Item($newRow, 15).Value = 26000
$ws.Cells.Item($newRow, 16).Value = 26000
$ws.Cells.Item($newRow, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item($newRow, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($newRow, 19).Value = 2600
$ws.Cells.Item($newRow, 20).Value = 10
